# Box Plot Updates, Color Updates Main Figures
#
# Re-positions the small data-point label textboxes (tx12..tx27) that sit
# inside the unnamed group shape on slide 1. Each pair of labels
# (tx12/tx13, tx14/tx15, ...) belongs to one box-plot whisker/point and
# both move together. Left/Top are expressed in points (PowerPoint's COM
# unit); the values below were derived from the target EMU offsets
# (EMU = points * 12700) so that round-tripping reproduces the exact
# OOXML <a:off> values.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)
$grp = $s.Shapes.Item(2)

function Set-ItemPos($index, $left, $top) {
    $shp = $grp.GroupItems.Item($index)
    $shp.Left = $left
    $shp.Top = $top
}

Set-ItemPos 10 453.4345   353.4966    # tx12
Set-ItemPos 11 466.2383   381.4557    # tx13
Set-ItemPos 12 318.7629   390.0483    # tx14
Set-ItemPos 13 337.7268   417.8157    # tx15
Set-ItemPos 14 309.2639   334.7628    # tx16
Set-ItemPos 15 329.1948   362.5219    # tx17
Set-ItemPos 16 262.9059   291.09261   # tx18
Set-ItemPos 17 277.6061   315.4591    # tx19
Set-ItemPos 18 340.79293  293.2441    # tx20
Set-ItemPos 19 364.5124   317.6106    # tx21
Set-ItemPos 20 288.6846   246.9862    # tx22
Set-ItemPos 21 293.4193   271.3526    # tx23
Set-ItemPos 22 327.6017   199.37993   # tx24
Set-ItemPos 23 347.5241   227.34741   # tx25
Set-ItemPos 24 431.61583  193.4983    # tx26
Set-ItemPos 25 440.6268   217.8564    # tx27
